$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "...cette licence. Kyl" -> "...cette licence jusqu'à WRC Generations. Kyl"
$d.Content.Find.Execute(
    "cette licence. Kyl", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "cette licence jusqu’à WRC Generations. Kyl", 2)

# --- Change 2 -----------------------------------------------------------
# "...du prochain WRC." -> "...du prochain titre (MMO open-world) Test Drive Unlimited."
$d.Content.Find.Execute(
    "du prochain WRC.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "du prochain titre (MMO open-world) Test Drive Unlimited.", 2)
